$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows (Min / Max columns) with new values
$ws.Range("B2").Value = 5.3
$ws.Range("C2").Value = 12.7

$ws.Range("B3").Value = 5.4
$ws.Range("C3").Value = 11.1

$ws.Range("B4").Value = 0.8
$ws.Range("C4").Value = 1.4

# Row 5 used to hold "theta_threshold_range" (with values 13.5 / 95) and is removed
# entirely, shifting "pie_threshold_range" (previously row 6) up to become row 5.
$ws.Rows("5").Delete()

# Apply the new values for the now-shifted pie_threshold_range row
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 20

# Update the selection / active cell to match the new view state
$ws.Range("C3").Select()

# Adjust the window size recorded for the workbook view
$wb.Windows.Item(1).Width = 26025 / 15
$wb.Windows.Item(1).Height = 9690 / 15

# Configure page setup so a pageSetup element (with printer settings) is written
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
